$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Acolyte" to "Outcast" (row 6, column A)
$ws.Range("A6").Value = "Outcast"
